$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 1638800
$ws.Range("E8").Value = 1829500
$ws.Range("F8").Value = 1317300
$ws.Range("G8").Value = 1133400
$ws.Range("H8").Value = 1007900
$ws.Range("I8").Value = 854800
$ws.Range("J8").Value = 702200
$ws.Range("D9").Value = 1565500
$ws.Range("E9").Value = 1639800
$ws.Range("F9").Value = 1278700
$ws.Range("G9").Value = 945500
$ws.Range("H9").Value = 789800
$ws.Range("I9").Value = 739600
$ws.Range("J9").Value = 659900
$ws.Range("D10").Value = 73300
$ws.Range("E10").Value = 189700
$ws.Range("F10").Value = 38600
$ws.Range("G10").Value = 187900
$ws.Range("H10").Value = 218100
$ws.Range("I10").Value = 115300
$ws.Range("J10").Value = 42400
$ws.Range("I12").Value = 700
$ws.Range("D14").Value = 23800
$ws.Range("E14").Value = 160900
$ws.Range("F14").Value = 229300
$ws.Range("G14").Value = -97600
$ws.Range("H14").Value = 66600
$ws.Range("I14").Value = 6800
$ws.Range("J14").Value = 2100
$ws.Range("D17").Value = 1598500
$ws.Range("E17").Value = 1808600
$ws.Range("F17").Value = 1516000
$ws.Range("G17").Value = 872500
$ws.Range("H17").Value = 865800
$ws.Range("I17").Value = 756000
$ws.Range("J17").Value = 672800
$ws.Range("D18").Value = 40300
$ws.Range("E18").Value = 20900
$ws.Range("F18").Value = -198700
$ws.Range("G18").Value = 260900
$ws.Range("H18").Value = 142100
$ws.Range("I18").Value = 98900
$ws.Range("J18").Value = 29400
$ws.Range("D20").Value = 38200
$ws.Range("E20").Value = -88800
$ws.Range("F20").Value = -78600
$ws.Range("G20").Value = -10700
$ws.Range("H20").Value = -101400
$ws.Range("I20").Value = -41600
$ws.Range("J20").Value = -23900
$ws.Range("D21").Value = 290600
$ws.Range("E21").Value = -19500
$ws.Range("F21").Value = -106200
$ws.Range("G21").Value = 260600
$ws.Range("H21").Value = 174000
$ws.Range("I21").Value = 86000
$ws.Range("J21").Value = 115700
$ws.Range("D22").Value = 67300
$ws.Range("E22").Value = 80400
$ws.Range("F22").Value = 80300
$ws.Range("D23").Value = 11200
$ws.Range("E23").Value = -148200
$ws.Range("F23").Value = -357600
$ws.Range("G23").Value = 250300
$ws.Range("H23").Value = 40700
$ws.Range("I23").Value = 57300
$ws.Range("J23").Value = 5500
$ws.Range("D24").Value = 5800
$ws.Range("E24").Value = -173600
$ws.Range("F24").Value = -28400
$ws.Range("G24").Value = 47800
$ws.Range("H24").Value = 34600
$ws.Range("I24").Value = 26200
$ws.Range("D26").Value = 5400
$ws.Range("E26").Value = 25300
$ws.Range("F26").Value = -329200
$ws.Range("G26").Value = 202500
$ws.Range("H26").Value = 6000
$ws.Range("I26").Value = 31100
$ws.Range("J26").Value = 5800
$ws.Range("D27").Value = 5300
$ws.Range("E27").Value = 25100
$ws.Range("F27").Value = -329200
$ws.Range("G27").Value = 215200
$ws.Range("H27").Value = 22800
$ws.Range("I27").Value = 36800
$ws.Range("J27").Value = 12300
$ws.Range("D32").Value = -38200
$ws.Range("E32").Value = 88800
$ws.Range("F32").Value = 78600
$ws.Range("G32").Value = 10700
$ws.Range("H32").Value = 101400
$ws.Range("I32").Value = 41600
$ws.Range("J32").Value = 23900
$ws.Range("D33").Value = 5300
$ws.Range("E33").Value = 25100
$ws.Range("F33").Value = -329200
$ws.Range("G33").Value = 215200
$ws.Range("H33").Value = 22800
$ws.Range("I33").Value = 36800
$ws.Range("J33").Value = 12300
$ws.Range("D35").Value = 5300
$ws.Range("E35").Value = 25100
$ws.Range("F35").Value = -329200
$ws.Range("G35").Value = 215200
$ws.Range("H35").Value = 22800
$ws.Range("I35").Value = 36800
$ws.Range("J35").Value = 12300
$ws.Range("D41").Value = 143900
$ws.Range("E41").Value = 141400
$ws.Range("F41").Value = 447100
$ws.Range("G41").Value = 66300
$ws.Range("H41").Value = 59700
$ws.Range("I41").Value = 49200
$ws.Range("J41").Value = 58600
$ws.Range("D43").Value = 448600
$ws.Range("E43").Value = 439700
$ws.Range("F43").Value = 400500
$ws.Range("G43").Value = 415200
$ws.Range("H43").Value = 146300
$ws.Range("I43").Value = 111500
$ws.Range("J43").Value = 71500
$ws.Range("D44").Value = 290300
$ws.Range("E44").Value = 241700
$ws.Range("F44").Value = 203000
$ws.Range("G44").Value = 46400
$ws.Range("H44").Value = 37800
$ws.Range("I44").Value = 27800
$ws.Range("J44").Value = 26500
$ws.Range("D45").Value = 87000
$ws.Range("D46").Value = 969900
$ws.Range("E46").Value = 822800
$ws.Range("F46").Value = 1050600
$ws.Range("G46").Value = 528000
$ws.Range("H46").Value = 243700
$ws.Range("I46").Value = 188500
$ws.Range("J46").Value = 156600
$ws.Range("D47").Value = 203900
$ws.Range("E47").Value = 173300
$ws.Range("F47").Value = 179500
$ws.Range("G47").Value = 172200
$ws.Range("H47").Value = 157900
$ws.Range("I47").Value = 11600
$ws.Range("J47").Value = 6800
$ws.Range("D48").Value = 3646700
$ws.Range("E48").Value = 3526000
$ws.Range("F48").Value = 3727100
$ws.Range("G48").Value = 1867100
$ws.Range("H48").Value = 1590900
$ws.Range("I48").Value = 1517000
$ws.Range("J48").Value = 1552300
$ws.Range("D49").Value = 480600
$ws.Range("E49").Value = 438400
$ws.Range("F49").Value = 466800
$ws.Range("G49").Value = 64200
$ws.Range("H49").Value = 66000
$ws.Range("I49").Value = 50500
$ws.Range("J49").Value = 50500
$ws.Range("D52").Value = 256500
$ws.Range("E52").Value = 253500
$ws.Range("F52").Value = 233700
$ws.Range("G52").Value = 228100
$ws.Range("H52").Value = 191500
$ws.Range("I52").Value = 169800
$ws.Range("J52").Value = 158700
$ws.Range("D54").Value = 5557500
$ws.Range("E54").Value = 5214000
$ws.Range("F54").Value = 5657600
$ws.Range("G54").Value = 2859600
$ws.Range("H54").Value = 2250000
$ws.Range("I54").Value = 1937300
$ws.Range("J54").Value = 1924900
$ws.Range("D57").Value = 468400
$ws.Range("E57").Value = 458600
$ws.Range("F57").Value = 490300
$ws.Range("G57").Value = 355100
$ws.Range("H57").Value = 266100
$ws.Range("I57").Value = 189100
$ws.Range("J57").Value = 196100
$ws.Range("D58").Value = 22900
$ws.Range("E58").Value = 113600
$ws.Range("F58").Value = 519000
$ws.Range("G58").Value = 51600
$ws.Range("H58").Value = 259100
$ws.Range("I58").Value = 136800
$ws.Range("J58").Value = 37200
$ws.Range("D59").Value = 33000
$ws.Range("E59").Value = 6200
$ws.Range("F59").Value = 9600
$ws.Range("G59").Value = 22200
$ws.Range("H59").Value = 34700
$ws.Range("I59").Value = 40600
$ws.Range("J59").Value = 38100
$ws.Range("D60").Value = 524300
$ws.Range("E60").Value = 578300
$ws.Range("F60").Value = 1018800
$ws.Range("G60").Value = 428800
$ws.Range("H60").Value = 559900
$ws.Range("I60").Value = 366500
$ws.Range("J60").Value = 271400
$ws.Range("D61").Value = 1943700
$ws.Range("E61").Value = 1644400
$ws.Range("F61").Value = 1551500
$ws.Range("G61").Value = 563500
$ws.Range("H61").Value = 222700
$ws.Range("I61").Value = 123900
$ws.Range("J61").Value = 217500
$ws.Range("D62").Value = 1340900
$ws.Range("E62").Value = 1346400
$ws.Range("F62").Value = 1437800
$ws.Range("G62").Value = 738500
$ws.Range("H62").Value = 483700
$ws.Range("I62").Value = 419800
$ws.Range("J62").Value = 437800
$ws.Range("D66").Value = 3810300
$ws.Range("E66").Value = 3570500
$ws.Range("F66").Value = 4009400
$ws.Range("G66").Value = 1732000
$ws.Range("H66").Value = 1256400
$ws.Range("I66").Value = 917800
$ws.Range("J66").Value = 941500
$ws.Range("D72").Value = -628900
$ws.Range("E72").Value = -732600
$ws.Range("F72").Value = -729900
$ws.Range("G72").Value = -362100
$ws.Range("H72").Value = -496100
$ws.Range("I72").Value = -470200
$ws.Range("J72").Value = -506200
$ws.Range("D76").Value = 1747200
$ws.Range("E76").Value = 1643500
$ws.Range("F76").Value = 1648200
$ws.Range("G76").Value = 1127600
$ws.Range("H76").Value = 993500
$ws.Range("I76").Value = 1019500
$ws.Range("J76").Value = 983500
$ws.Range("D81").Value = 5300
$ws.Range("E81").Value = 25100
$ws.Range("F81").Value = -329200
$ws.Range("G81").Value = 215200
$ws.Range("H81").Value = 22800
$ws.Range("I81").Value = 36800
$ws.Range("J81").Value = 12300
$ws.Range("D89").Value = 160000
$ws.Range("E89").Value = 113000
$ws.Range("F89").Value = 74900
$ws.Range("G89").Value = 132500
$ws.Range("H89").Value = 169500
$ws.Range("I89").Value = 146300
$ws.Range("J89").Value = 94600
$ws.Range("D91").Value = -210100
$ws.Range("E91").Value = -247800
$ws.Range("F91").Value = -170300
$ws.Range("G91").Value = -163800
$ws.Range("H91").Value = -120700
$ws.Range("I91").Value = -122500
$ws.Range("J91").Value = -106700
$ws.Range("D94").Value = -207600
$ws.Range("E94").Value = -253200
$ws.Range("F94").Value = -1675800
$ws.Range("G94").Value = -268100
$ws.Range("H94").Value = -379200
$ws.Range("I94").Value = -123400
$ws.Range("J94").Value = -105500
$ws.Range("D100").Value = 41200
$ws.Range("E100").Value = -151400
$ws.Range("F100").Value = 1988700
$ws.Range("G100").Value = 142200
$ws.Range("H100").Value = 231000
$ws.Range("I100").Value = -32300
$ws.Range("J100").Value = 30800
$ws.Range("D101").Value = 9000
$ws.Range("E101").Value = -14100
$ws.Range("F101").Value = -7100
$ws.Range("H101").Value = -10800
$ws.Range("D102").Value = 2500
$ws.Range("E102").Value = -305700
$ws.Range("F102").Value = 380800
$ws.Range("G102").Value = 6700
$ws.Range("H102").Value = 10500
$ws.Range("I102").Value = -9400
$ws.Range("J102").Value = 20000
